$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '91.130.80'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  +0.63%  '
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.156.57'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  +1.38%  '
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '239.57'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  +0.64%  '
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '618.99'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -0.99%  '
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.12'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  +0.33%  '
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.374'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  +0.51%  '
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '3.155.13'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  +1.37%  '
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.736'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  +0.26%  '
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.203'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  +0.22%  '
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000247'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -0.80%  '
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.64'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  +3.31%  '
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '34.99'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -0.74%  '
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '90.820.93'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  +0.68%  '
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.746.30'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  +1.65%  '
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.161.72'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  +1.93%  '
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -3.81%  '
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.02'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  +5.39%  '
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.95'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  +3.97%  '
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '457.52'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  +2.73%  '
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.0000205'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -3.72%  '
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.15'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  +0.31%  '
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '6.00'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  +1.27%  '
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '88.90'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -1.63%  '
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.94'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -1.38%  '
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = 'Hedera'
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.153'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  +43.08%  '
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = 'WrappedeETH'
$ws.Cells.Item(29, 2).Style = "Normal"
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.333.80'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  +5.35%  '
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -0.07%  '
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.236'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  +8.27%  '
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = 'Cronos'
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.168'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -3.28%  '
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 2).Style = "Normal"
$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '9.49'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  +2.52%  '
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.173'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +11.65%  '
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -0.99%  '
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '26.58'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -0.13%  '
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '7.62'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  +4.15%  '
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.96'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  +1.98%  '
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '501.80'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  +1.16%  '
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.33'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  +3.53%  '
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.86'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -8.74%  '
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  +7.07%  '
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.50'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -2.79%  '
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '22.12'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  +0.11%  '
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  +0.01%  '
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.715'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  +4.94%  '
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  +1.58%  '
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '155.78'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -1.56%  '
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.37'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  +3.19%  '
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '4.47'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -1.22%  '
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = 'OKB'
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(51, 3).Style = "Normal"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '44.18'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -1.48%  '
$ws.Cells.Item(51, 5).Style = "Normal"
